# ----------------------------------------------------------------------
# This script applies the "feat: add 2022-Q3 data" edit:
#  1. Insert a new row into the "总计" (summary) sheet for 2022-Q3.
#  2. Insert a brand-new worksheet "2022-Q3" (placed right after "总计")
#     containing the per-fund holdings data for that quarter.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sum = $wb.Worksheets.Item(1)

# --- Part 1: "总计" sheet -> insert new row 2 for 2022-Q3 -----------------
$sum.Rows.Item(2).Insert()
$sum.Range("B2:D2").ClearFormats()
$sum.Range("A2").Value = 0
$sum.Range("B2").Value = "2022-Q3"
$sum.Range("C2").Value = 18
$sum.Range("D2").Value = 3.26
# Row-Insert does not carry the bordered/bold/centered style (s=2) that
# column A uses on every data row, so copy it back from the row below.
$sum.Range("A3").Copy()
$sum.Range("A2").PasteSpecial(-4122)
$sum.Range("A2").Value = 0

# --- Part 2: add the new "2022-Q3" worksheet ------------------------------
$newSheet = $wb.Worksheets.Add($null, $sum)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $newSheet.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Each row: @(A-index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名)
$data = @(
    @(0, "009951", "广发稳健回报混合A", "52.44", "46.49", "3.00", "1.5732", 3),
    @(1, "005505", "前海开源中药研究精选股票A", "4.10", "83.96", "7.13", "0.2923", 8),
    @(2, "009952", "广发稳健回报混合C", "9.48", "46.49", "3.00", "0.2844", 3),
    @(3, "005506", "前海开源中药研究精选股票C", "3.58", "83.96", "7.13", "0.2553", 8),
    @(4, "501070", "广发睿阳三年定期开放混合", "6.62", "51.01", "3.14", "0.2079", 9),
    @(5, "005775", "中加转型动力灵活配置混合A", "5.34", "50.55", "2.62", "0.1399", 6),
    @(6, "673110", "西部利得新润灵活配置混合A", "4.48", "76.88", "2.84", "0.1272", 6),
    @(7, "010457", "广发睿鑫混合A", "3.08", "73.79", "3.68", "0.1133", 7),
    @(8, "163001", "长信医疗保健行业灵活配置混合（LOF）", "1.88", "94.04", "4.71", "0.0885", 9),
    @(9, "009242", "中加核心智造混合A", "1.92", "61.20", "2.47", "0.0474", 9),
    @(10, "007254", "广发均衡价值混合", "0.66", "89.74", "6.60", "0.0436", 2),
    @(11, "005776", "中加转型动力灵活配置混合C", "1.27", "50.55", "2.62", "0.0333", 6),
    @(12, "010458", "广发睿鑫混合C", "0.53", "73.79", "3.68", "0.0195", 7),
    @(13, "002872", "华夏智胜价值成长股票C", "2.13", "93.39", "0.85", "0.0181", 8),
    @(14, "002871", "华夏智胜价值成长股票A", "0.86", "93.39", "0.85", "0.0073", 8),
    @(15, "013154", "长信医疗保健行业灵活配置混合(LOF)C", "0.05", "94.04", "4.71", "0.0024", 9),
    @(16, "009243", "中加核心智造混合C", "0.09", "61.20", "2.47", "0.0022", 9),
    @(17, "015356", "西部利得新润灵活配置混合C", "0.00", "76.88", "2.84", "0", 6),
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    # Prefix text-like numeric strings with an apostrophe so Excel keeps them
    # as text instead of silently converting them to numbers.
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($r -eq 19) {
        # last row holds position value 0 as a real number, not text
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the throw-away "text" number-format/style created by the apostrophe
# trick above, then re-apply the real styles used throughout the workbook:
#  - header row (row 1, columns B:H) uses the bold/bordered/centered style
#  - column A on every data row uses the same bold/bordered/centered style
$newSheet.Range("B1:H1").ClearFormats()
$newSheet.Range("B2:G19").ClearFormats()
$sum.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$sum.Range("A2").Copy()
$newSheet.Range("A2:A19").PasteSpecial(-4122)

